$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Absent") is derived from column E ("Real"): Absent = 1 - Real.
# Recompute H for every data row (rows 3-21) to form the consolidated report.
for ($r = 3; $r -le 21; $r++) {
    $real = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 8).Value = 1 - $real
}
